{"js": "// Update the date paragraph and all of the answer cells in the practice\n// table. Every text run in the document changes, so we address the date\n// paragraph directly and each table cell by (row, col) to avoid any\n// ambiguity from duplicate cell text (e.g. \"36\u00f79=4, 0\" appears twice with\n// different replacements).\n\nconst body = context.document.body;\n\n// --- Date heading -----------------------------------------------------\n// The heading is the document's first paragraph (everything else lives\n// inside the table), so address it positionally rather than searching for\n// the old date text.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.getRange().insertText(\"2025-06-04 Wednesday\", \"Replace\");\n\n// --- Table of division answers -----------------------------------------\nconst table = body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\n// New values, keyed by row index within the 5x5 answer grid (rows 0, 4, 8,\n// 12, 16 of the 20-row table hold text; the others are spacer rows).\nconst newValues = [\n  [\"26\u00f79=2, 8\", \"15\u00f78=1, 7\", \"60\u00f76=10, 0\", \"42\u00f78=5, 2\", \"53\u00f79=5, 8\"],\n  [\"44\u00f77=6, 2\", \"44\u00f73=14, 2\", \"88\u00f73=29, 1\", \"71\u00f73=23, 2\", \"94\u00f73=31, 1\"],\n  [\"52\u00f72=26, 0\", \"13\u00f74=3, 1\", \"87\u00f72=43, 1\", \"67\u00f73=22, 1\", \"59\u00f73=19, 2\"],\n  [\"21\u00f79=2, 3\", \"78\u00f75=15, 3\", \"57\u00f76=9, 3\", \"93\u00f79=10, 3\", \"86\u00f77=12, 2\"],\n  [\"26\u00f79=2, 8\", \"64\u00f78=8, 0\", \"52\u00f75=10, 2\", \"62\u00f73=20, 2\", \"60\u00f76=10, 0\"],\n];\n\nconst answerRows = [0, 4, 8, 12, 16];\n\nfor (let i = 0; i < answerRows.length; i++) {\n  const row = answerRows[i];\n  for (let col = 0; col < newValues[i].length; col++) {\n    table.getCell(row, col).value = newValues[i][col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date paragraph and every answer cell in the practice table.\n# Every text run in the document changes, so cells are addressed directly\n# by (row, col) via the Tables/Cell object model rather than a text search\n# -- several original cell values repeat (e.g. \"36\u00f79=4, 0\" occurs twice)\n# but map to different replacements, so a blind Find/Replace would be\n# ambiguous.\n\n$d = $word.ActiveDocument\n\n# --- Date heading -----------------------------------------------------\n$d.Paragraphs.Item(1).Range.Text = \"2025-06-04 Wednesday\"\n\n# --- Table of division answers -----------------------------------------\n$tbl = $d.Tables.Item(1)\n\n# New values, keyed by the 1-based table row that holds text (rows 1, 5, 9,\n# 13, 17 of the 20-row table; the rows between them are blank spacers).\n$answerRows = 1, 5, 9, 13, 17\n\n$newValues = @(\n  @(\"26\u00f79=2, 8\", \"15\u00f78=1, 7\", \"60\u00f76=10, 0\", \"42\u00f78=5, 2\", \"53\u00f79=5, 8\"),\n  @(\"44\u00f77=6, 2\", \"44\u00f73=14, 2\", \"88\u00f73=29, 1\", \"71\u00f73=23, 2\", \"94\u00f73=31, 1\"),\n  @(\"52\u00f72=26, 0\", \"13\u00f74=3, 1\", \"87\u00f72=43, 1\", \"67\u00f73=22, 1\", \"59\u00f73=19, 2\"),\n  @(\"21\u00f79=2, 3\", \"78\u00f75=15, 3\", \"57\u00f76=9, 3\", \"93\u00f79=10, 3\", \"86\u00f77=12, 2\"),\n  @(\"26\u00f79=2, 8\", \"64\u00f78=8, 0\", \"52\u00f75=10, 2\", \"62\u00f73=20, 2\", \"60\u00f76=10, 0\")\n)\n\nfor ($i = 0; $i -lt $answerRows.Length; $i++) {\n  $row = $answerRows[$i]\n  for ($col = 1; $col -le 5; $col++) {\n    $tbl.Cell($row, $col).Range.Text = $newValues[$i][$col - 1]\n  }\n}\n"}
